$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.24

$ws.Range("D3").Value = 1.41
$ws.Range("E3").Value = 1.33

$ws.Range("B4").Value = 1.5
$ws.Range("F4").Value = 1.1

$ws.Range("C5").Value = 1.34
$ws.Range("F5").Value = 1.06

$ws.Range("G6").Value = 0.98

$ws.Range("F7").Value = 1.49
